# "fixes to distored hierarchy"
#
# The rows for class_29 (B106:B126) on sheet "3_import" were mis-classified;
# they belong under class_30 (matching the rows immediately below them).
# Fix the hierarchy values, then leave the workbook's view state positioned
# the way the author left it when saving: "3_import" active/selected with
# B106:B126 highlighted (the block that was just corrected), and "2_import"
# no longer the selected tab.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("2_import")
$ws3 = $wb.Worksheets.Item("3_import")

# Correct the distorted hierarchy: class_29 -> class_30 for rows 106-126.
$ws3.Range("B106:B126").Value = "class_30"

# Move the active tab / selection to "3_import", highlighting the fixed
# range, and make "2_import" no longer the selected sheet.
$ws3.Activate() | Out-Null
$ws3.Range("B106:B126").Select() | Out-Null
$ws2.Range("D21").Select() | Out-Null
$ws3.Activate() | Out-Null
